$wb = $excel.ActiveWorkbook

# "Logs" sheet gets a new row (21) appended with the latest mail-log entry.
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A21").Value = "Demo inplannen"
$logs.Range("B21").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C21").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D21").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E21").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F21").Value = "2025-08-14 21:13:48"
$logs.Range("G21").Value = "Nee"
$logs.Range("H21").Value = "Ja"
$logs.Range("I21").Value = "Nee"
$logs.Range("J21").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too,
# same as Excel does automatically when the underlying table/range grows.
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))
$logs.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J21"))

# "Dashboard" sheet: bump the count for the category that just got a new row.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 15
